$d = $word.ActiveDocument

$replacements = @(
    @{old="322×8="; new="841×7="},
    @{old="495×4="; new="362×4="},
    @{old="358×6="; new="711×2="},
    @{old="381×8="; new="191×9="},
    @{old="789×3="; new="273×5="},
    @{old="738×4="; new="411×4="},
    @{old="156×2="; new="245×5="},
    @{old="426×5="; new="448×9="},
    @{old="587×2="; new="437×5="},
    @{old="140×2="; new="490×3="},
    @{old="516×2="; new="683×2="},
    @{old="485×3="; new="522×8="},
    @{old="848×4="; new="149×2="},
    @{old="259×9="; new="169×8="},
    @{old="630×9="; new="888×6="},
    @{old="529×2="; new="261×3="},
    @{old="451×6="; new="836×6="},
    @{old="540×4="; new="480×2="},
    @{old="681×3="; new="750×3="},
    @{old="774×8="; new="344×9="},
    @{old="487×7="; new="163×2="},
    @{old="390×8="; new="935×7="},
    @{old="972×4="; new="874×9="},
    @{old="949×5="; new="945×6="},
    @{old="496×7="; new="885×2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
